# Update header labels on Sheet1 to be capitalized, matching the target
# shared-strings content (a -> A, b -> B, c -> C, d -> D, message -> Message).
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("B1").Value = "A"
$ws1.Range("C1").Value = "B"
$ws1.Range("D1").Value = "C"
$ws1.Range("E1").Value = "D"
$ws1.Range("F1").Value = "Message"

# Remove Sheet2 entirely.
$excel.DisplayAlerts = $false
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Delete()
$excel.DisplayAlerts = $true
